$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - "Save", matching the style of the other header cells (e.g. G1)
# PasteSpecial with xlPasteFormats (-4122) only copies formatting (bold font,
# thin border, centered/top alignment), so it reuses the existing style index
# instead of Excel minting a brand-new (duplicate) style record.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data values for column H (Save), rows 2-8
$saveValues = @(0, 1, 1, 0, 0, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
